$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (M2:T2)
$ws.Range("M2").Value = 0.8317113333333332
$ws.Range("N2").Value = 2.495134
$ws.Range("O2").Value = 0.0263454906755698
$ws.Range("P2").Value = 0.0263454906755698
$ws.Range("Q2").Value = 1.064169383494889
$ws.Range("R2").Value = 9.577524451453998
$ws.Range("S2").Value = 0.0263454906755698
$ws.Range("T2").Value = 0.0263454906755698

# Row 3 updates (O3, P3, S3, T3)
$ws.Range("O3").Value = 0.6529848313028861
$ws.Range("P3").Value = 0.6529848313028862
$ws.Range("S3").Value = 0.6529848313028861
$ws.Range("T3").Value = 0.6529848313028862

# Row 4 updates (M4:T4)
$ws.Range("M4").Value = 10.12334933333333
$ws.Range("N4").Value = 30.370048
$ws.Range("O4").Value = 0.3206696780215441
$ws.Range("P4").Value = 0.3206696780215441
$ws.Range("Q4").Value = 12.95276135745422
$ws.Range("R4").Value = 116.574852217088
$ws.Range("S4").Value = 0.3206696780215441
$ws.Range("T4").Value = 0.3206696780215441
